$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -7
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = -7
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 7
$ws.Range("F12").Value = -4
